# Burnin Shift Activities - roll the shift date from 2025-02-17 to 2025-02-18,
# update the engineer sign-off to "Administrador", swap in the new "GRR EN LA
# SECCION 6" activity in row 3, and log the previous "EJECUTAR GRR" activity
# as a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new shift date.
$ws.Name = "Shift 2 - 2025-02-18"

# Column A holds the shift date for every data row - bump it from the 17th to
# the 18th everywhere it appears. A leading apostrophe keeps Excel from
# re-interpreting the ISO-looking text as a serial date.
$ws.Range("A2").Value = "'2025-02-18"
$ws.Range("A3").Value = "'2025-02-18"
$ws.Range("A4").Value = "'2025-02-18"

# Row 3 now documents the GRR run against section 6 (with its description),
# replacing what used to be the "EJECUTAR GRR" / "GRR 2H-23 2H-38" entry.
$ws.Range("C3").Value = "GRR EN LA SECCION 6"
$ws.Range("D3").Value = "Esta si tiene descripcion"

# The sign-off on row 4 is shortened to just "Administrador".
$ws.Range("E4").Value = "Administrador"

# Log the activity that used to live in row 3 as a brand-new row 5, signed
# off the same way as row 4.
$ws.Range("A5").Value = "'2025-02-18"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "EJECUTAR GRR"
$ws.Range("D5").Value = "GRR 2H-23 2H-38"
$ws.Range("E5").Value = "Administrador"

# Match the visual style already used by the other data rows (thin border all
# round, text wrapped, vertically top-aligned).
$dataRng = $ws.Range("A5:E5")
$dataRng.Borders.LineStyle = 1
$dataRng.WrapText = $true
$dataRng.VerticalAlignment = -4160
